$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 36524
$ws.Range("D2").Value = 52859587
$ws.Range("C3").Value = 88578
$ws.Range("D3").Value = 129943082
$ws.Range("C4").Value = 30354
$ws.Range("D4").Value = 44979848
$ws.Range("C5").Value = 8410
$ws.Range("D5").Value = 12503176
$ws.Range("C6").Value = 1859
$ws.Range("D6").Value = 2764025
$ws.Range("C7").Value = 139
$ws.Range("D7").Value = 203593
$ws.Range("C11").Value = 39951
$ws.Range("D11").Value = 54283306
$ws.Range("C12").Value = 9368
$ws.Range("D12").Value = 13552253
$ws.Range("C13").Value = 25372
$ws.Range("D13").Value = 37229074
$ws.Range("C14").Value = 8128
$ws.Range("D14").Value = 12064817
$ws.Range("C15").Value = 2090
$ws.Range("D15").Value = 3110150
$ws.Range("C16").Value = 394
$ws.Range("D16").Value = 580123
$ws.Range("C19").Value = 9872
$ws.Range("D19").Value = 13109585
$ws.Range("C20").Value = 13090
$ws.Range("D20").Value = 18906844
$ws.Range("C21").Value = 31050
$ws.Range("D21").Value = 45591790
$ws.Range("C22").Value = 10053
$ws.Range("D22").Value = 14950497
$ws.Range("C23").Value = 2562
$ws.Range("D23").Value = 3813263
$ws.Range("C24").Value = 472
$ws.Range("D24").Value = 701839
$ws.Range("C25").Value = 35
$ws.Range("D25").Value = 51953
$ws.Range("C26").Value = 11395
$ws.Range("D26").Value = 15249066
$ws.Range("C27").Value = 7431
$ws.Range("D27").Value = 10764358
$ws.Range("C28").Value = 22047
$ws.Range("D28").Value = 32362408
$ws.Range("C29").Value = 7624
$ws.Range("D29").Value = 11345095
$ws.Range("C30").Value = 1920
$ws.Range("D30").Value = 2865544
$ws.Range("C31").Value = 347
$ws.Range("D31").Value = 517915
$ws.Range("C33").Value = 8112
$ws.Range("D33").Value = 10724758
$ws.Range("C34").Value = 3089
$ws.Range("D34").Value = 4455809
$ws.Range("C35").Value = 7485
$ws.Range("D35").Value = 10937337
$ws.Range("C36").Value = 3018
$ws.Range("D36").Value = 4471288
$ws.Range("C37").Value = 792
$ws.Range("D37").Value = 1180763
$ws.Range("C38").Value = 145
$ws.Range("D38").Value = 215732
$ws.Range("C40").Value = 2316
$ws.Range("D40").Value = 3127817
$ws.Range("C41").Value = 16773
$ws.Range("D41").Value = 24273210
$ws.Range("C42").Value = 49844
$ws.Range("D42").Value = 73104465
$ws.Range("C43").Value = 18578
$ws.Range("D43").Value = 27597353
$ws.Range("C44").Value = 5429
$ws.Range("D44").Value = 8085222
$ws.Range("C45").Value = 1107
$ws.Range("D45").Value = 1651292
$ws.Range("C49").Value = 16232
$ws.Range("D49").Value = 21643474
$ws.Range("C50").Value = 1889
$ws.Range("D50").Value = 2740705
$ws.Range("C51").Value = 6511
$ws.Range("D51").Value = 9577473
$ws.Range("C52").Value = 2253
$ws.Range("D52").Value = 3364690
$ws.Range("C53").Value = 734
$ws.Range("D53").Value = 1096305
$ws.Range("C54").Value = 171
$ws.Range("D54").Value = 253333
$ws.Range("C56").Value = 6191
$ws.Range("D56").Value = 8521785
$ws.Range("C57").Value = 842
$ws.Range("D57").Value = 1233743
$ws.Range("C58").Value = 2036
$ws.Range("D58").Value = 3021826
$ws.Range("C59").Value = 833
$ws.Range("D59").Value = 1240337
$ws.Range("C60").Value = 291
$ws.Range("D60").Value = 436258
$ws.Range("C61").Value = 73
$ws.Range("D61").Value = 109500
$ws.Range("C63").Value = 1209
$ws.Range("D63").Value = 1707230
$ws.Range("C64").Value = 14977
$ws.Range("D64").Value = 21640831
$ws.Range("C65").Value = 43774
$ws.Range("D65").Value = 64090776
$ws.Range("C66").Value = 15386
$ws.Range("D66").Value = 22875406
$ws.Range("C67").Value = 4460
$ws.Range("D67").Value = 6644060
$ws.Range("C68").Value = 893
$ws.Range("D68").Value = 1329096
$ws.Range("C71").Value = 14730
$ws.Range("D71").Value = 19447554
$ws.Range("C72").Value = 49204
$ws.Range("D72").Value = 71635217
$ws.Range("C73").Value = 140818
$ws.Range("D73").Value = 207545311
$ws.Range("C74").Value = 61480
$ws.Range("D74").Value = 91637667
$ws.Range("C75").Value = 19564
$ws.Range("D75").Value = 29233959
$ws.Range("C76").Value = 4493
$ws.Range("D76").Value = 6713472
$ws.Range("C83").Value = 48875
$ws.Range("D83").Value = 66659810
$ws.Range("C84").Value = 4415
$ws.Range("D84").Value = 6398185
$ws.Range("C85").Value = 11158
$ws.Range("D85").Value = 16398798
$ws.Range("C86").Value = 3779
$ws.Range("D86").Value = 5631915
$ws.Range("C87").Value = 1316
$ws.Range("D87").Value = 1966489
$ws.Range("C91").Value = 5151
$ws.Range("D91").Value = 6933713
$ws.Range("C92").Value = 1495
$ws.Range("D92").Value = 2161041
$ws.Range("C93").Value = 4917
$ws.Range("D93").Value = 7243114
$ws.Range("C94").Value = 1873
$ws.Range("D94").Value = 2790999
$ws.Range("C95").Value = 660
$ws.Range("D95").Value = 988960
$ws.Range("C99").Value = 3338
$ws.Range("D99").Value = 4432876
$ws.Range("C100").Value = 566
$ws.Range("D100").Value = 844464
$ws.Range("C101").Value = 326
$ws.Range("D101").Value = 487130
$ws.Range("C102").Value = 114
$ws.Range("D102").Value = 171000
$ws.Range("C105").Value = 10478
$ws.Range("D105").Value = 15222890
$ws.Range("C106").Value = 28640
$ws.Range("D106").Value = 42087346
$ws.Range("C107").Value = 9592
$ws.Range("D107").Value = 14264257
$ws.Range("C108").Value = 2622
$ws.Range("D108").Value = 3909307
$ws.Range("C109").Value = 465
$ws.Range("D109").Value = 694982
$ws.Range("C112").Value = 9550
$ws.Range("D112").Value = 12640360
$ws.Range("C113").Value = 29581
$ws.Range("D113").Value = 42688026
$ws.Range("C114").Value = 64670
$ws.Range("D114").Value = 94676689
$ws.Range("C115").Value = 20917
$ws.Range("D115").Value = 31097479
$ws.Range("C116").Value = 5889
$ws.Range("D116").Value = 8774780
$ws.Range("C117").Value = 1070
$ws.Range("D117").Value = 1599506
$ws.Range("C118").Value = 69
$ws.Range("D118").Value = 100920
$ws.Range("C121").Value = 25064
$ws.Range("D121").Value = 33512898
$ws.Range("C122").Value = 34698
$ws.Range("D122").Value = 50113192
$ws.Range("C123").Value = 74582
$ws.Range("D123").Value = 109134392
$ws.Range("C124").Value = 23211
$ws.Range("D124").Value = 34459905
$ws.Range("C125").Value = 6177
$ws.Range("D125").Value = 9180443
$ws.Range("C126").Value = 1153
$ws.Range("D126").Value = 1713911
$ws.Range("C130").Value = 30600
$ws.Range("D130").Value = 40699591
$ws.Range("C131").Value = 12904
$ws.Range("D131").Value = 18688016
$ws.Range("C132").Value = 31662
$ws.Range("D132").Value = 46528380
$ws.Range("C133").Value = 11265
$ws.Range("D133").Value = 16739059
$ws.Range("C134").Value = 2889
$ws.Range("D134").Value = 4308081
$ws.Range("C135").Value = 466
$ws.Range("D135").Value = 692990
$ws.Range("C138").Value = 10546
$ws.Range("D138").Value = 14079301
$ws.Range("C139").Value = 33941
$ws.Range("D139").Value = 49042135
$ws.Range("C140").Value = 79048
$ws.Range("D140").Value = 115852890
$ws.Range("C141").Value = 23735
$ws.Range("D141").Value = 35285365
$ws.Range("C142").Value = 6194
$ws.Range("D142").Value = 9245586
$ws.Range("C143").Value = 1365
$ws.Range("D143").Value = 2029785
$ws.Range("C146").Value = 28316
$ws.Range("D146").Value = 38283212
